# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '69.374.75'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = '3.687.10'
$ws.Range("E3").Value = '  -0.06%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '680.29'
$ws.Range("E5").Value = '  -1.63%  '

$ws.Range("D6").Value = '159.37'
$ws.Range("E6").Value = '  -2.16%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  -1.10%  '

$ws.Range("E9").Value = '  -1.42%  '

$ws.Range("D10").Value = '7.05'
$ws.Range("E10").Value = '  -4.55%  '

$ws.Range("E11").Value = '  -1.84%  '

$ws.Range("E12").Value = '  -3.47%  '

$ws.Range("D13").Value = '4.311.51'
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").Value = '32.43'
$ws.Range("E14").Value = '  -3.34%  '

$ws.Range("D15").Value = '3.695.55'
$ws.Range("E15").Value = '  +0.39%  '

$ws.Range("D16").Value = '69.329.09'
$ws.Range("E16").Value = '  -0.12%  '

$ws.Range("E17").Value = '  +1.75%  '

$ws.Range("D18").Value = '16.04'
$ws.Range("E18").Value = '  -1.73%  '

$ws.Range("D19").Value = '6.42'
$ws.Range("E19").Value = '  -2.99%  '

$ws.Range("D20").Value = '468.46'
$ws.Range("E20").Value = '  -2.86%  '

$ws.Range("D21").Value = '10.02'
$ws.Range("E21").Value = '  -0.06%  '

$ws.Range("E22").Value = '  -2.41%  '

$ws.Range("D23").Value = '79.94'
$ws.Range("E23").Value = '  -0.16%  '

$ws.Range("D24").Value = '3.834.30'
$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("E26").Value = '  -5.72%  '

$ws.Range("E27").Value = '  -4.47%  '

$ws.Range("D28").Value = '9.13'
$ws.Range("E28").Value = '  -4.62%  '

$ws.Range("E29").Value = '  -1.95%  '

$ws.Range("D30").Value = '1.76'
$ws.Range("E30").Value = '  -3.79%  '

$ws.Range("E31").Value = '  -3.65%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '1.99'
$ws.Range("E32").Value = '  -4.68%  '

$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.06%  '

$ws.Range("D34").Value = '26.94'
$ws.Range("E34").Value = '  -0.61%  '

$ws.Range("D35").Value = '3.676.94'
$ws.Range("E35").Value = '  +0.63%  '

$ws.Range("E36").Value = '  -4.89%  '

$ws.Range("D37").Value = '8.28'
$ws.Range("E37").Value = '  -3.16%  '

$ws.Range("D38").Value = '6.25'
$ws.Range("E38").Value = '  -1.63%  '

$ws.Range("D40").Value = '2.27'
$ws.Range("E40").Value = '  -2.58%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("E42").Value = '  -3.14%  '

$ws.Range("D43").Value = '171.65'
$ws.Range("E43").Value = '  +4.80%  '

$ws.Range("E44").Value = '  -1.04%  '

$ws.Range("D45").Value = '47.72'
$ws.Range("E45").Value = '  -0.70%  '

$ws.Range("E46").Value = '  -4.04%  '

$ws.Range("D47").Value = '28.12'
$ws.Range("E47").Value = '  -6.64%  '

$ws.Range("D48").Value = '0.000277'
$ws.Range("E48").Value = '  -3.44%  '

$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Value = '1.29'
$ws.Range("E49").Value = '  -4.60%  '

$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").Value = '1.10'
$ws.Range("E50").Value = '  -4.02%  '

$ws.Range("E51").Value = '  -2.91%  '
